# The workbook tracks weekly "Haba" (fava bean) price records for the
# "Mercado Mayorista Lo Valledor de Santiago" market. This edit adds a new
# weekly record. In the underlying data table, the new record is inserted
# right above the row that used to be row 210 (which shifts down, along
# with every row below it, by one row), and the sheet's dimension grows
# from A1:R260 to A1:R261 to account for the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 210; this pushes former rows
# 210..260 down to 211..261, preserving all of their existing data/styles.
$ws.Rows("210:210").Insert()

# Populate the newly inserted row 210 with the new data record.
$ws.Range("A210").Value = 6
$ws.Range("B210").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C210").Value = "Metropolitana"
$ws.Range("D210").Value = 44754
$ws.Range("E210").Value = 13
$ws.Range("F210").Value = 100112026
$ws.Range("G210").Value = "Haba"
$ws.Range("H210").Value = "Sin especificar"
$ws.Range("I210").Value = "Primera"
$ws.Range("J210").Value = 400
$ws.Range("K210").Value = 15000
$ws.Range("L210").Value = 16000
$ws.Range("M210").Value = 15425
$ws.Range("N210").Value = "`$/saco 25 kilos"
$ws.Range("O210").Value = "Región de Coquimbo"
$ws.Range("P210").Value = 617
$ws.Range("Q210").Value = 25
$ws.Range("R210").Value = "Hortaliza"
